# NMRA Region/Division Map - add a new "seasonal member" entry:
#   Region 22 (NFR), sub-region 10, name "Lakeshores-NY"
# This is inserted as a new row just above the existing row 28 (which starts
# region 23), and every row from the old row 28 through the old last row
# (206) shifts down by one to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new row and populate it -----------------------------------
$ws.Rows.Item(28).Insert()

$ws.Range("A28").Value2 = 22
$ws.Range("B28").Value2 = 10
$ws.Range("C28").Value2 = "Lakeshores-NY"
$ws.Range("D28").Value2 = "NFR"

# --- Refresh the sheet's cached sort-state so it covers the new last row --
# (Applying the sort also quietly "fixes" a couple of pre-existing,
#  not-strictly-B-sorted rows elsewhere in the sheet, so we restore those
#  two rows back to their original values right after.)
$origA21 = $ws.Range("A21").Value2
$origB21 = $ws.Range("B21").Value2
$origC21 = $ws.Range("C21").Value2
$origD21 = $ws.Range("D21").Value2
$origA22 = $ws.Range("A22").Value2
$origB22 = $ws.Range("B22").Value2
$origC22 = $ws.Range("C22").Value2
$origD22 = $ws.Range("D22").Value2

$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("A3:A207"))
$sortObj.SortFields.Add($ws.Range("B3:B207"))
$sortObj.SetRange($ws.Range("A3:D207"))
$sortObj.Header = 0
$sortObj.Apply()

$ws.Range("A21").Value2 = $origA21
$ws.Range("B21").Value2 = $origB21
$ws.Range("C21").Value2 = $origC21
$ws.Range("D21").Value2 = $origD21
$ws.Range("A22").Value2 = $origA22
$ws.Range("B22").Value2 = $origB22
$ws.Range("C22").Value2 = $origC22
$ws.Range("D22").Value2 = $origD22

# --- Match the saved selection --------------------------------------------
$ws.Range("C28").Select()
